# PI-2023-24 Self-Assessment workbook update
# "Added the rest of the team photos" -- fills in the remaining peer
# review rows on the "User Stories" sheet, bumps a handful of rubric
# scores on "Project Management", consolidates the per-column
# conditional-formatting rules on "User Stories" into a single rule
# covering the whole grid, and widens/repositions the related data
# validation ranges + view selections to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: User Stories  (peer grades + newly added team members)
# ---------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("User Stories")

# Existing reviewers: self-assessment grade bumped 4 -> 5
$wsUser.Range("C6").Value  = 5
$wsUser.Range("C7").Value  = 5
$wsUser.Range("C8").Value  = 5
$wsUser.Range("C11").Value = 5
$wsUser.Range("C12").Value = 5
$wsUser.Range("C13").Value = 5

# Newly added team members (rows 14-19), each with student #, ID and grade
$wsUser.Range("A14").Value = 9
$wsUser.Range("B14").Value = 1230481
$wsUser.Range("C14").Value = 4

$wsUser.Range("A15").Value = 10
$wsUser.Range("B15").Value = 1230929
$wsUser.Range("C15").Value = 4

$wsUser.Range("A16").Value = 11
$wsUser.Range("B16").Value = 1231151
$wsUser.Range("C16").Value = 4

$wsUser.Range("A17").Value = 12
$wsUser.Range("B17").Value = 1231170
$wsUser.Range("C17").Value = 4

$wsUser.Range("A18").Value = 13
$wsUser.Range("B18").Value = 1231151
$wsUser.Range("C18").Value = 4

$wsUser.Range("A19").Value = 14
$wsUser.Range("B19").Value = 1231151
$wsUser.Range("C19").Value = 4

# Consolidate the 17 per-column conditional formats into a single rule
# spanning the whole grades grid.
$wsUser.Cells.FormatConditions.Delete()
$fc = $wsUser.Range("E6:J25").FormatConditions.Add(2, 3, "`$C6=E`$3")
$fc.StopIfTrue = $true
$fc.Font.Color = 24832
$fc.Interior.Color = 13561798

# Data validation ranges grow to cover the newly filled-in rows.
$wsUser.Range("C18:C25").Validation.Delete()
$wsUser.Range("C20:C25").Validation.Add(3, 1, 1, "=`$E`$40:`$J`$40")
$wsUser.Range("C6:C17").Validation.Delete()
$wsUser.Range("C6:C19").Validation.Add(3, 1, 1, "=`$E`$3:`$J`$3")

# ---------------------------------------------------------------------
# Sheet: Project Management (rubric score touch-ups)
# ---------------------------------------------------------------------
$wsPM = $wb.Worksheets.Item("Project Management")

foreach ($col in @("C", "D", "E", "F", "G")) {
    $wsPM.Range($col + "4").Value = 4
}

foreach ($col in @("C", "D", "E", "F", "G")) {
    $wsPM.Range($col + "8").Value = 5
}

$wsPM.Range("C9").Value = 4
$wsPM.Range("D9").Value = 4
$wsPM.Range("E9").Value = 4
$wsPM.Range("F9").Value = 5
$wsPM.Range("G9").Value = 4

foreach ($col in @("C", "D", "E", "F", "G")) {
    $wsPM.Range($col + "10").Value = 4
}

foreach ($col in @("C", "D", "E", "F", "G")) {
    $wsPM.Range($col + "11").Value = 3
}

# ---------------------------------------------------------------------
# Sheet view / selection bookkeeping (matches the saved cursor state)
# ---------------------------------------------------------------------
$wsGroup = $wb.Worksheets.Item("Group and Self Assessment")
$wsGroup.Select()
$wsGroup.Range("J15").Select()

$wsUser.Select()
$excel.ActiveWindow.Zoom = 61
$wsUser.Range("B18").Select()

$wsProjDev = $wb.Worksheets.Item("Project Development")
$wsProjDev.Select()
$wsProjDev.Range("I5").Select()

$wsPM.Select()
$wsPM.Range("I10").Select()
